$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = 0.09
$ws.Range("B10").Value = 23.065251918
$ws.Range("C10").Value = 19.391128632
$ws.Range("D10").Value = 3.674123286
$ws.Range("E10").Value = 15.134183304
$ws.Range("F10").Value = 19.013517068
$ws.Range("G10").Value = 16.23480563
$ws.Range("H10").Value = 2.778711438
$ws.Range("I10").Value = 16.287015978

$ws.Range("A11").Value = 0.1
$ws.Range("B11").Value = 53.477439452
$ws.Range("C11").Value = 27.743483948
$ws.Range("D11").Value = 25.733955504
$ws.Range("E11").Value = 19.02164175333333
$ws.Range("F11").Value = 33.43191751
$ws.Range("G11").Value = 20.299145172
$ws.Range("H11").Value = 13.132772338
$ws.Range("I11").Value = 20.22482003466667

$ws.Range("A12").Value = 0.11
$ws.Range("B12").Value = 762.840795484
$ws.Range("C12").Value = 35.700369212
$ws.Range("D12").Value = 727.1404262719999
$ws.Range("E12").Value = 22.97236418
$ws.Range("F12").Value = 357.411734862
$ws.Range("G12").Value = 23.752440764
$ws.Range("H12").Value = 333.659294098
$ws.Range("I12").Value = 24.18630784933334

$ws.Range("A13").Value = 0.12
$ws.Range("B13").Value = 1953.745469976
$ws.Range("C13").Value = 38.77228725200001
$ws.Range("D13").Value = 1914.973182722
$ws.Range("E13").Value = 24.811487422
$ws.Range("F13").Value = 872.5247670800001
$ws.Range("G13").Value = 24.700685818
$ws.Range("H13").Value = 847.824081262
$ws.Range("I13").Value = 26.029398638

$ws.Range("A14").Value = 0.13
$ws.Range("B14").Value = 2266.191827468
$ws.Range("C14").Value = 39.018818624
$ws.Range("D14").Value = 2227.173008844
$ws.Range("E14").Value = 25.243859218
$ws.Range("F14").Value = 986.171924806
$ws.Range("G14").Value = 24.503048226
$ws.Range("H14").Value = 961.6688765800001
$ws.Range("I14").Value = 26.46850563866667

$ws.Range("A15").Value = 0.14
$ws.Range("B15").Value = 3032.198468768
$ws.Range("C15").Value = 38.950103506
$ws.Range("D15").Value = 2993.248365262
$ws.Range("E15").Value = 25.52585647666667
$ws.Range("F15").Value = 1284.978745966
$ws.Range("G15").Value = 24.161410788
$ws.Range("H15").Value = 1260.817335178
$ws.Range("I15").Value = 26.75116729333333

$ws.Range("A16").Value = 0.15
$ws.Range("B16").Value = 3467.50966688
$ws.Range("C16").Value = 38.289160512
$ws.Range("D16").Value = 3429.220506368
$ws.Range("E16").Value = 25.51369319333333
$ws.Range("F16").Value = 1436.474800022
$ws.Range("G16").Value = 23.65727812
$ws.Range("H16").Value = 1412.817521902
$ws.Range("I16").Value = 26.749201902

$ws.Range("A17").Value = 0.16
$ws.Range("B17").Value = 3911.892387694
$ws.Range("C17").Value = 37.487058824
$ws.Range("D17").Value = 3874.405328872
$ws.Range("E17").Value = 25.440260226
$ws.Range("F17").Value = 1585.105019658
$ws.Range("G17").Value = 23.140807668
$ws.Range("H17").Value = 1561.96421199
$ws.Range("I17").Value = 26.67706112333333

$ws.Range("A18").Value = 0.17
$ws.Range("B18").Value = 4041.011614314
$ws.Range("C18").Value = 36.964767414
$ws.Range("D18").Value = 4004.0468469
$ws.Range("E18").Value = 25.44409872933333
$ws.Range("F18").Value = 1608.56998527
$ws.Range("G18").Value = 22.774044308
$ws.Range("H18").Value = 1585.795940964
$ws.Range("I18").Value = 26.68560294533333

$ws.Range("A19").Value = 0.18
$ws.Range("B19").Value = 4312.789545273999
$ws.Range("C19").Value = 36.30853267600001
$ws.Range("D19").Value = 4276.481012597999
$ws.Range("E19").Value = 25.30718837866667
$ws.Range("F19").Value = 1696.902866872
$ws.Range("G19").Value = 22.475940104
$ws.Range("H19").Value = 1674.426926766
$ws.Range("I19").Value = 26.55429984466667

$ws.Range("A20").Value = 0.19
$ws.Range("B20").Value = 4427.371668214
$ws.Range("C20").Value = 35.672811882
$ws.Range("D20").Value = 4391.698856334
$ws.Range("E20").Value = 25.31023158866667
$ws.Range("F20").Value = 1708.110787172
$ws.Range("G20").Value = 22.116897024
$ws.Range("H20").Value = 1685.993890148
$ws.Range("I20").Value = 26.56685283933334

$ws.Range("A21").Value = 0.2
$ws.Range("B21").Value = 4719.73852585
$ws.Range("C21").Value = 35.270357802
$ws.Range("D21").Value = 4684.468168048
$ws.Range("E21").Value = 25.32196972466667
$ws.Range("F21").Value = 1798.145568474
$ws.Range("G21").Value = 21.926013394
$ws.Range("H21").Value = 1776.21955508
$ws.Range("I21").Value = 26.58711745933334
